# Apply scheduled-runner profit recalculation updates across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2137.2104
$ws.Range("I96").Value = 593.3570999999999
$ws.Range("J96").Value = 6460
$ws.Range("K96").Value = 1780.0713
$ws.Range("L96").Value = 19380
$ws.Range("M96").Value = -407.0712999999998
$ws.Range("N96").Value = -22126
$ws.Range("H98").Value = 5055743.5
$ws.Range("I98").Value = 5546.4
$ws.Range("J98").Value = 55557716
$ws.Range("K98").Value = 5546.4
$ws.Range("L98").Value = 55557716
$ws.Range("M98").Value = -4048.4
$ws.Range("N98").Value = -55560712
$ws.Range("H100").Value = 14494782
$ws.Range("I100").Value = 19608874
$ws.Range("K100").Value = 19608874
$ws.Range("M100").Value = -19608333
$ws.Range("H112").Value = 3152.5938
$ws.Range("J112").Value = 3378.862
$ws.Range("L112").Value = 10136.586
$ws.Range("N112").Value = -12352.586
$ws.Range("H122").Value = 5055743.5
$ws.Range("I122").Value = 5546.4
$ws.Range("J122").Value = 55557716
$ws.Range("K122").Value = 16639.2
$ws.Range("L122").Value = 166673148
$ws.Range("M122").Value = -14189.2
$ws.Range("N122").Value = -166678048
$ws.Range("H125").Value = 2150
$ws.Range("I125").Value = 2100
$ws.Range("J125").Value = 2200
$ws.Range("K125").Value = 18900
$ws.Range("L125").Value = 19800
$ws.Range("M125").Value = -16440
$ws.Range("N125").Value = -24720
$ws.Range("H137").Value = 1340.0526
$ws.Range("I137").Value = 960.1818
$ws.Range("K137").Value = 2880.5454
$ws.Range("M137").Value = -330.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5834.5713
$ws.Range("I74").Value = 7452.4
$ws.Range("J74").Value = 1790
$ws.Range("K74").Value = 7452.4
$ws.Range("L74").Value = 1790
$ws.Range("M74").Value = -6578.4
$ws.Range("N74").Value = -3538
$ws.Range("H77").Value = 5834.5713
$ws.Range("I77").Value = 7452.4
$ws.Range("J77").Value = 1790
$ws.Range("K77").Value = 37262
$ws.Range("L77").Value = 8950
$ws.Range("M77").Value = -32894
$ws.Range("N77").Value = -17686

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 9500.75
$ws.Range("J51").Value = 9500.75
$ws.Range("L51").Value = 9500.75
$ws.Range("N51").Value = -10972.75
$ws.Range("H58").Value = 1826.9333
$ws.Range("J58").Value = 2068.2222
$ws.Range("L58").Value = 2068.2222
$ws.Range("N58").Value = -2474.2222
$ws.Range("H61").Value = 9500.75
$ws.Range("J61").Value = 9500.75
$ws.Range("L61").Value = 9500.75
$ws.Range("N61").Value = -10196.75
$ws.Range("H80").Value = 20975
$ws.Range("I80").Value = 10000
$ws.Range("K80").Value = 10000
$ws.Range("M80").Value = -8877
$ws.Range("H83").Value = 20975
$ws.Range("I83").Value = 10000
$ws.Range("K83").Value = 30000
$ws.Range("M83").Value = -24384
$ws.Range("H99").Value = 2355
$ws.Range("I99").Value = 2634.5
$ws.Range("J99").Value = 2106.5557
$ws.Range("K99").Value = 2634.5
$ws.Range("L99").Value = 2106.5557
$ws.Range("M99").Value = -1136.5
$ws.Range("N99").Value = -5102.5557
$ws.Range("H122").Value = 1834.8235
$ws.Range("I122").Value = 1479.6923
$ws.Range("J122").Value = 2989
$ws.Range("K122").Value = 4439.0769
$ws.Range("L122").Value = 8967
$ws.Range("M122").Value = -1989.0769
$ws.Range("N122").Value = -13867
$ws.Range("H126").Value = 2355
$ws.Range("I126").Value = 2634.5
$ws.Range("J126").Value = 2106.5557
$ws.Range("K126").Value = 7903.5
$ws.Range("L126").Value = 6319.6671
$ws.Range("M126").Value = -5433.5
$ws.Range("N126").Value = -11259.6671
$ws.Range("H134").Value = 1733.1333
$ws.Range("I134").Value = 1345.9231
$ws.Range("J134").Value = 4250
$ws.Range("K134").Value = 4037.7693
$ws.Range("L134").Value = 12750
$ws.Range("M134").Value = -1502.7693
$ws.Range("N134").Value = -17820
$ws.Range("H136").Value = 1826.9333
$ws.Range("J136").Value = 2068.2222
$ws.Range("L136").Value = 6204.6666
$ws.Range("N136").Value = -11304.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 459.8889
$ws.Range("I5").Value = 494.07693
$ws.Range("J5").Value = 371
$ws.Range("K5").Value = 1482.23079
$ws.Range("L5").Value = 1113
$ws.Range("M5").Value = -1370.23079
$ws.Range("N5").Value = -1337
$ws.Range("H92").Value = 600
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 0
$ws.Range("N92").Value = -552
$ws.Range("M92").ClearContents()
$ws.Range("H113").Value = 804.16
$ws.Range("I113").Value = 500.80768
$ws.Range("J113").Value = 910.7432
$ws.Range("K113").Value = 1502.42304
$ws.Range("L113").Value = 2732.2296
$ws.Range("M113").Value = 667.5769599999999
$ws.Range("N113").Value = -7072.229600000001
$ws.Range("H131").Value = 782.62
$ws.Range("J131").Value = 820.9425
$ws.Range("L131").Value = 2462.8275
$ws.Range("N131").Value = -12542.8275
$ws.Range("H135").Value = 459.8889
$ws.Range("I135").Value = 494.07693
$ws.Range("J135").Value = 371
$ws.Range("K135").Value = 4446.69237
$ws.Range("L135").Value = 3339
$ws.Range("M135").Value = -1911.69237
$ws.Range("N135").Value = -8409

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1390.6316
$ws.Range("I102").Value = 1241.8572
$ws.Range("J102").Value = 1807.2
$ws.Range("K102").Value = 1241.8572
$ws.Range("L102").Value = 1807.2
$ws.Range("M102").Value = 380.1428000000001
$ws.Range("N102").Value = -5051.2
$ws.Range("H122").Value = 2999.4
$ws.Range("J122").Value = 4500
$ws.Range("L122").Value = 13500
$ws.Range("N122").Value = -18400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2567.3572
$ws.Range("I7").Value = 2107.923
$ws.Range("J7").Value = 3313.9375
$ws.Range("K7").Value = 2107.923
$ws.Range("L7").Value = 3313.9375
$ws.Range("M7").Value = -1995.923
$ws.Range("N7").Value = -3537.9375
$ws.Range("H126").Value = 2567.3572
$ws.Range("I126").Value = 2107.923
$ws.Range("J126").Value = 3313.9375
$ws.Range("K126").Value = 6323.768999999999
$ws.Range("L126").Value = 9941.8125
$ws.Range("M126").Value = -3853.768999999999
$ws.Range("N126").Value = -14881.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2933.9
$ws.Range("I81").Value = 1948.25
$ws.Range("J81").Value = 4412.375
$ws.Range("K81").Value = 3896.5
$ws.Range("L81").Value = 8824.75
$ws.Range("M81").Value = -2835.5
$ws.Range("N81").Value = -10946.75
$ws.Range("H84").Value = 2933.9
$ws.Range("I84").Value = 1948.25
$ws.Range("J84").Value = 4412.375
$ws.Range("K84").Value = 19482.5
$ws.Range("L84").Value = 44123.75
$ws.Range("M84").Value = -14178.5
$ws.Range("N84").Value = -54731.75
$ws.Range("H126").Value = 1407.3334
$ws.Range("I126").Value = 1156.6666
$ws.Range("J126").Value = 1992.2222
$ws.Range("K126").Value = 3469.9998
$ws.Range("L126").Value = 5976.6666
$ws.Range("M126").Value = -999.9998000000001
$ws.Range("N126").Value = -10916.6666
